# Testing/Test Cases/Login TC .xlsx
# Commit: Added test case BANK_SYS_TC_Log_R025
#
# Row 32 previously held test case BANK_SYS_TC_Log_R024 (the "press back /
# press forward" scenario) but had a typo ("userenter") and empty
# priority/type columns. This fixes the typo and fills in the priority
# ("high") and type ("functional") columns.
#
# Row 33 was a blank placeholder row; it now becomes the new test case
# BANK_SYS_TC_Log_R025 (internet disconnects before home page loads after
# login).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 32: fix typo + fill in priority / type ----
$ws.Cells.Item(32, 2).Value = "BANK_SYS_TC_Log_R024"
$ws.Cells.Item(32, 3).Value = "verify that if user enter the 3 fields `"user name `",`"password `" ,`"National ID`"  then he pressed back before he press log in button then he pressed forward the login page appear with empty fields `n"
$ws.Cells.Item(32, 8).Value = "high"
$ws.Cells.Item(32, 9).Value = "functional"

# ---- Row 33: new test case BANK_SYS_TC_Log_R025 ----
$ws.Cells.Item(33, 2).Value = "BANK_SYS_TC_Log_R025"
$ws.Cells.Item(33, 3).Value = "verify that if user enter the 3 fields `"user name `",`"password `" ,`"National ID`"  then he pressed log in then the internet disconnected before the home page load `nafter the internet connection back`nhe login page appear with empty fields "
$ws.Cells.Item(33, 4).Value = "user name:`npassword:`nNational ID:"
$ws.Cells.Item(33, 5).Value = "1-enter   user name field  with  `"user name`" like in the test data`n2-enter password like in test data`n3-enter   National ID field  with like in the test`n4-press login and disconnected the internnt at the same time "
$ws.Cells.Item(33, 6).Value = "the login page appear with empty fields "
$ws.Cells.Item(33, 8).Value = "high"
$ws.Cells.Item(33, 9).Value = "functional"

# Row 33 grew taller to fit the new content.
$ws.Rows.Item(33).RowHeight = 135.75

# Move the view / selection down to the newly added row, like the author did.
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B33").Select()
